$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.039666533470154
$ws.Range("B1").Value = 1.659794688224792
$ws.Range("C1").Value = 6.846080303192139
$ws.Range("D1").Value = 2.701626062393188
$ws.Range("E1").Value = 1.478082656860352
